# A new weekly record was inserted into the price table at row 102,
# pushing the existing rows 102-233 down to 103-234 (dimension grows
# from A1:R233 to A1:R234).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 102, shifting everything
# below it (including formatting) down by one row.
$ws.Rows("102:102").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A102").Value = 5
$ws.Range("B102").Value = "Macroferia Regional de Talca"
$ws.Range("C102").Value = "Maule"
$ws.Range("D102").Value = 44671
$ws.Range("E102").Value = 7
$ws.Range("F102").Value = 100112009
$ws.Range("G102").Value = "Acelga"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 500
$ws.Range("K102").Value = 3500
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = 3500
$ws.Range("N102").Value = "$/docena de atados (4 kilos)"
$ws.Range("O102").Value = "Región del Maule"
$ws.Range("P102").Value = 875
$ws.Range("Q102").Value = 4
$ws.Range("R102").Value = "Hortaliza"
